# The workbook has two sheets: "2021-Q3" (the big per-fund holdings table,
# A1:H32) and "总计" (the small summary table, A1:D2). The edit re-sorts the
# sheet order/names so that "总计" becomes the first sheet and "2021-Q3" the
# second, while each sheet's data "follows" its new name/position (i.e. the
# physical sheet that ends up named "总计" carries the small summary table,
# and the one named "2021-Q3" carries the big holdings table).

$wb = $excel.ActiveWorkbook

# Add a scratch sheet to hold both tables temporarily while we cross-swap
# them. Worksheets.Add() inserts at index 1, which pushes the two real
# sheets to index 2 / 3 -- fetch fresh references to them afterwards.
$tmp = $wb.Worksheets.Add()
$shA = $wb.Worksheets.Item(2)   # currently "2021-Q3" -- the big A1:H32 table
$shB = $wb.Worksheets.Item(3)   # currently "总计"     -- the small A1:D2 table

# Stash both tables (values + formatting) into the scratch sheet, far apart
# so the source ranges never overlap the destinations.
$shA.Range("A1:H32").Copy($tmp.Range("A1"))
$shB.Range("A1:D2").Copy($tmp.Range("Z1"))

# Wipe both real sheets completely so nothing stray is left behind.
$shA.Cells.Clear() | Out-Null
$shB.Cells.Clear() | Out-Null

# Cross-paste: the sheet that used to be "2021-Q3" gets the small summary
# table, and the sheet that used to be "总计" gets the big holdings table.
$tmp.Range("Z1:AC2").Copy($shA.Range("A1"))
$tmp.Range("A1:H32").Copy($shB.Range("A1"))

# Drop the scratch sheet; indices shift back to 1 / 2.
$tmp.Delete() | Out-Null
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ws1 now holds the "总计" data, ws2 now holds the "2021-Q3" data -- rename
# to match. Route through a scratch name so the two renames never collide
# on a duplicate sheet name mid-swap.
$ws1.Name = "__SWAP_TMP__"
$ws2.Name = "2021-Q3"
$ws1.Name = "总计"

# "2021-Q3" (now the second sheet) is the selected/active tab.
$ws2.Activate()
